# ---------------------------------------------------------------------------
# indicadores.xlsx – update several KPI figures on sheet "ÁREA GESTIÓN
# CORPORATIVA". The "CONSOLIDADO ACADÉMICO" sheet pulls these via formulas,
# so it recalculates automatically once the source cells change.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ÁREA GESTIÓN CORPORATIVA")

# E3: 10+15+4+5 (=34)  ->  10+15+4+5-5 (=29)
$ws.Range("E3").Formula = "=10+15+4+5-5"

# E4: 12+6 (=18)  ->  12+6+5 (=23)
$ws.Range("E4").Formula = "=12+6+5"

# F5: 16/43 (≈0.3721)  ->  21/43 (≈0.4884)
$ws.Range("F5").Formula = "=21/43"

# C6: 0.0284  ->  0.0919
$ws.Range("C6").Value = 0.0919

# F6: 0.0303  ->  0.097
$ws.Range("F6").Value = 0.097

# F20: 7  ->  9
$ws.Range("F20").Value = 9

# Force a full recalculation so the dependent "CONSOLIDADO ACADÉMICO"
# formulas (B3, C4, C5) pick up the new figures.
$excel.CalculateFullRebuild()

# Reflect the new selection / scroll position recorded for this sheet in the
# saved workbook (activeCell moved from G23 to C11, view scrolled to B1).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C11").Select()
